$d = $word.ActiveDocument

$pairs = @(
    @("203÷6=", "522÷2="),
    @("865÷9=", "964÷6="),
    @("297÷8=", "280÷4="),
    @("332÷8=", "931÷3="),
    @("460÷6=", "722÷5="),
    @("904÷9=", "730÷4="),
    @("109÷5=", "143÷9="),
    @("635÷8=", "668÷2="),
    @("429÷8=", "720÷5="),
    @("504÷6=", "565÷2="),
    @("822÷8=", "632÷6="),
    @("846÷9=", "148÷2="),
    @("208÷5=", "105÷5="),
    @("507÷5=", "999÷9="),
    @("648÷8=", "763÷6="),
    @("354÷4=", "302÷2="),
    @("590÷5=", "524÷5="),
    @("864÷4=", "762÷5="),
    @("910÷2=", "955÷4="),
    @("323÷3=", "887÷2="),
    @("246÷6=", "946÷6="),
    @("678÷5=", "366÷7="),
    @("150÷7=", "125÷6="),
    @("196÷5=", "308÷3="),
    @("686÷9=", "900÷5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
